# Daily attendance processing - 2025-11-09 18:26:49
# Normalizes the "Recorded By" (column G) entries: the literal "System"
# token is moved out of the leading position in the comma-separated list
# (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"),
# except when "System" is paired only with "backup@backdoor.com", whose
# order is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) { continue }

    $parts = $value -split ", "

    if ($parts.Length -eq 2 -and $parts[0].Equals("System") -and -not $parts[1].Equals("backup@backdoor.com")) {
        $cell.Value = $parts[1] + ", " + $parts[0]
    }
    elseif ($parts.Length -eq 3 -and $parts[0].Equals("System") -and $parts[1].Equals("system") -and $parts[2].Equals("backup@backdoor.com")) {
        $cell.Value = $parts[1] + ", " + $parts[0] + ", " + $parts[2]
    }
}
